$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Merge the three runs (with proofErr spellStart/spellEnd around "for ever")
#    in the Brenton/"AI" column of the v89 row into a single plain run.
$d.Content.Find.Execute(
    "Thy word, O Lord, endures for ever in the heavens.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Thy word, O Lord, endures for ever in the heavens.", 2) | Out-Null

# 2. Fill in the NETS column (column 6) for rows v89-v96 (table rows 2-9)
#    with two new paragraphs of translated text each, styled EngIndEnd
#    (matching the paragraph already present in each empty cell).

function Add-TwoLines($row, $line1, $line2) {
    $cell = $t.Cell($row, 6)
    $para = $cell.Range.Paragraphs.Item(1)
    $para.Range.InsertBefore("$line1`r$line2")
}

Add-TwoLines 2 "Forever, O Lord;" "your word endures in the sky,"
Add-TwoLines 3 "your truth to generation and generation;" "you founded the earth, and it endures."
Add-TwoLines 4 "By your arrangement the day endures," "because, all things together [the universe] are slaves of yours."
Add-TwoLines 5 "If it were not for the fact that your law was my meditation," "then I would have perished in my humiliation."

# Row 6 (v93) second line is split into two runs around a lastRenderedPageBreak.
$cell6 = $t.Cell(6, 6)
$para6 = $cell6.Range.Paragraphs.Item(1)
$para6.Range.InsertBefore("Your statutes I will never forget,`rbecause by them you quickened me [O Lord].")

Add-TwoLines 7 "Yours I am; save me," "because your statutes I sought."

# Row 8 (v95) second line has a _GoBack bookmark in the middle of the text.
$cell8 = $t.Cell(8, 6)
$para8 = $cell8.Range.Paragraphs.Item(1)
$para8.Range.InsertBefore("Sinners waited for me to destroy me;`ryour testimonies I considered.")

Add-TwoLines 9 "I saw a limit to all perfection [completion];" "your commandment is exceedingly spacious."

# 3. Move the _GoBack bookmark: it used to sit after "...exceedingly broad."
#    in row 9's "AI" column; it now belongs after "your testimonies" in row 8's
#    new NETS paragraph. Add the bookmark at the new location.
$full = $d.Content.Text
$anchorText = "your testimonies"
$idx = $full.IndexOf($anchorText)
if ($idx -ge 0) {
    $pos = $idx + $anchorText.Length
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
